$d = $word.ActiveDocument

# 1. Title
$d.Content.Find.Execute(
    "Hands-on AI-Assisted Programming Made Simple with GitHub Copilot",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Github Foundations Certification Training", 2)

# 2. WSQ funding sentence (course name already replaced by step 1)
$d.Content.Find.Execute(
    "according to Digital Technology Adoption and Innovation ACC-ICT-3004-1.1 under Infocomm Technology Framework.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "according to Software Configuration ICT-DIT-3014-1.1 under Infocomm Technology Framework.", 2)

# 3. Performance gap paragraph
$d.Content.Find.Execute(
    "One significant challenge is the slow adoption of new technologies and methodologies, hindering the ability to remain competitive. Teams may lack the expertise to effectively integrate emerging tools, which can lead to missed opportunities for automation and improved productivity. Legacy systems and a reluctance to change further compound these issues.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Integrating and deploying software products often involves using a variety of scripts and tools, but many developers lack the expertise to select the most appropriate ones. This can lead to inefficient deployment processes, increased errors, and difficulty in maintaining software quality across different platforms and environments. This impacts developer productivity and the overall reliability of software deployments.", 2)

# 4. Course description paragraph
$d.Content.Find.Execute(
    "This course directly addresses this by providing hands-on experience with cutting-edge AI programming tools. Participants will explore how these tools can streamline organizational coding processes, including using code completion and suggesting code snippets. The course is designed to help individuals stay current with the latest technology and propose relevant IT solutions.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This course teaches participants to use specific scripts and tools tailored for integrating and deploying software, greatly enhancing their skill set. The course emphasizes choosing the right tools for the job, which leads to more efficient deployments and better overall integration, which reduces deployment issues and improves developer output.", 2)

# 5. Date
$d.Content.Find.Execute(
    ": 03 March 2025",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": 04 March 2025", 2)

Write-Output "done"
